# Apply updated NATMI ligand-receptor statistics (C3-Nrp1 sheet)
# per "Natmi following Dr Hou advice" - recompute with updated
# ligand/receptor expressing-cell counts (1 -> 3) and dependent stats.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.558564
$ws.Range("H2").Value = 4.675692
$ws.Range("I2").Value = 0.005692101168584756
$ws.Range("J2").Value = 0.005692101168584756
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 118.0470123333333
$ws.Range("N2").Value = 354.141037
$ws.Range("O2").Value = 0.4657216250363638
$ws.Range("P2").Value = 0.4657216250363638
$ws.Range("Q2").Value = 183.9838237302893
$ws.Range("R2").Value = 1655.854413572604
$ws.Range("S2").Value = 0.002650934606104678
$ws.Range("T2").Value = 0.002650934606104678

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.558564
$ws.Range("H3").Value = 4.675692
$ws.Range("I3").Value = 0.005692101168584756
$ws.Range("J3").Value = 0.005692101168584756
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 59.36586533333332
$ws.Range("N3").Value = 178.097596
$ws.Range("O3").Value = 0.2342114953037475
$ws.Range("P3").Value = 0.2342114953037476
$ws.Range("Q3").Value = 92.5255005373813
$ws.Range("R3").Value = 832.7295048364318
$ws.Range("S3").Value = 0.001333155526114444
$ws.Range("T3").Value = 0.001333155526114445

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.558564
$ws.Range("H4").Value = 4.675692
$ws.Range("I4").Value = 0.005692101168584756
$ws.Range("J4").Value = 0.005692101168584756
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 56.84506433333333
$ws.Range("N4").Value = 170.535193
$ws.Range("O4").Value = 0.2242663767030476
$ws.Range("P4").Value = 0.2242663767030477
$ws.Range("Q4").Value = 88.59667084761732
$ws.Range("R4").Value = 797.3700376285559
$ws.Range("S4").Value = 0.001276546904905687
$ws.Range("T4").Value = 0.001276546904905687

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.558564
$ws.Range("H5").Value = 4.675692
$ws.Range("I5").Value = 0.005692101168584756
$ws.Range("J5").Value = 0.005692101168584756
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 19.21324333333333
$ws.Range("N5").Value = 57.63973
$ws.Range("O5").Value = 0.07580050295684103
$ws.Range("P5").Value = 0.07580050295684104
$ws.Range("Q5").Value = 29.94506938257333
$ws.Range("R5").Value = 269.50562444316
$ws.Range("S5").Value = 0.0004314641314599471
$ws.Range("T5").Value = 0.0004314641314599472

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 259.5505726666667
$ws.Range("H6").Value = 778.6517180000001
$ws.Range("I6").Value = 0.9479162344201305
$ws.Range("J6").Value = 0.9479162344201304
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 118.0470123333333
$ws.Range("N6").Value = 354.141037
$ws.Range("O6").Value = 0.4657216250363638
$ws.Range("P6").Value = 0.4657216250363638
$ws.Range("Q6").Value = 30639.16965270573
$ws.Range("R6").Value = 275752.5268743516
$ws.Range("S6").Value = 0.4414650890924939
$ws.Range("T6").Value = 0.4414650890924939

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 259.5505726666667
$ws.Range("H7").Value = 778.6517180000001
$ws.Range("I7").Value = 0.9479162344201305
$ws.Range("J7").Value = 0.9479162344201304
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 59.36586533333332
$ws.Range("N7").Value = 178.097596
$ws.Range("O7").Value = 0.2342114953037475
$ws.Range("P7").Value = 0.2342114953037476
$ws.Range("Q7").Value = 15408.44434411888
$ws.Range("R7").Value = 138675.9990970699
$ws.Range("S7").Value = 0.2220128786862365
$ws.Range("T7").Value = 0.2220128786862365

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 259.5505726666667
$ws.Range("H8").Value = 778.6517180000001
$ws.Range("I8").Value = 0.9479162344201305
$ws.Range("J8").Value = 0.9479162344201304
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 56.84506433333333
$ws.Range("N8").Value = 170.535193
$ws.Range("O8").Value = 0.2242663767030476
$ws.Range("P8").Value = 0.2242663767030477
$ws.Range("Q8").Value = 14754.16900099018
$ws.Range("R8").Value = 132787.5210089116
$ws.Range("S8").Value = 0.2125857393113994
$ws.Range("T8").Value = 0.2125857393113994

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 259.5505726666667
$ws.Range("H9").Value = 778.6517180000001
$ws.Range("I9").Value = 0.9479162344201305
$ws.Range("J9").Value = 0.9479162344201304
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 19.21324333333333
$ws.Range("N9").Value = 57.63973
$ws.Range("O9").Value = 0.07580050295684103
$ws.Range("P9").Value = 0.07580050295684104
$ws.Range("Q9").Value = 4986.808309950683
$ws.Range("R9").Value = 44881.27478955615
$ws.Range("S9").Value = 0.07185252733000072
$ws.Range("T9").Value = 0.07185252733000072

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 12.18925266666667
$ws.Range("H10").Value = 36.567758
$ws.Range("I10").Value = 0.04451691386950307
$ws.Range("J10").Value = 0.04451691386950307
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 118.0470123333333
$ws.Range("N10").Value = 354.141037
$ws.Range("O10").Value = 0.4657216250363638
$ws.Range("P10").Value = 0.4657216250363638
$ws.Range("Q10").Value = 1438.904859876116
$ws.Range("R10").Value = 12950.14373888504
$ws.Range("S10").Value = 0.02073248946890881
$ws.Range("T10").Value = 0.02073248946890881

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 12.18925266666667
$ws.Range("H11").Value = 36.567758
$ws.Range("I11").Value = 0.04451691386950307
$ws.Range("J11").Value = 0.04451691386950307
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 59.36586533333332
$ws.Range("N11").Value = 178.097596
$ws.Range("O11").Value = 0.2342114953037475
$ws.Range("P11").Value = 0.2342114953037476
$ws.Range("Q11").Value = 723.6255323233074
$ws.Range("R11").Value = 6512.629790909767
$ws.Range("S11").Value = 0.01042637296368445
$ws.Range("T11").Value = 0.01042637296368445

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 12.18925266666667
$ws.Range("H12").Value = 36.567758
$ws.Range("I12").Value = 0.04451691386950307
$ws.Range("J12").Value = 0.04451691386950307
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 56.84506433333333
$ws.Range("N12").Value = 170.535193
$ws.Range("O12").Value = 0.2242663767030476
$ws.Range("P12").Value = 0.2242663767030477
$ws.Range("Q12").Value = 692.8988520119216
$ws.Range("R12").Value = 6236.089668107294
$ws.Range("S12").Value = 0.009983646975515102
$ws.Range("T12").Value = 0.009983646975515103

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 12.18925266666667
$ws.Range("H13").Value = 36.567758
$ws.Range("I13").Value = 0.04451691386950307
$ws.Range("J13").Value = 0.04451691386950307
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 19.21324333333333
$ws.Range("N13").Value = 57.63973
$ws.Range("O13").Value = 0.07580050295684103
$ws.Range("P13").Value = 0.07580050295684104
$ws.Range("Q13").Value = 234.1950775361489
$ws.Range("R13").Value = 2107.75569782534
$ws.Range("S13").Value = 0.003374404461394705
$ws.Range("T13").Value = 0.003374404461394706

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.5133286666666667
$ws.Range("H14").Value = 1.539986
$ws.Range("I14").Value = 0.001874750541781658
$ws.Range("J14").Value = 0.001874750541781658
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 118.0470123333333
$ws.Range("N14").Value = 354.141037
$ws.Range("O14").Value = 0.4657216250363638
$ws.Range("P14").Value = 0.4657216250363638
$ws.Range("Q14").Value = 60.59691544505355
$ws.Range("R14").Value = 545.3722390054819
$ws.Range("S14").Value = 0.0008731118688563572
$ws.Range("T14").Value = 0.0008731118688563572

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.5133286666666667
$ws.Range("H15").Value = 1.539986
$ws.Range("I15").Value = 0.001874750541781658
$ws.Range("J15").Value = 0.001874750541781658
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 59.36586533333332
$ws.Range("N15").Value = 178.097596
$ws.Range("O15").Value = 0.2342114953037475
$ws.Range("P15").Value = 0.2342114953037476
$ws.Range("Q15").Value = 30.47420049707288
$ws.Range("R15").Value = 274.2678044736559
$ws.Range("S15").Value = 0.000439088127712193
$ws.Range("T15").Value = 0.0004390881277121931

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.5133286666666667
$ws.Range("H16").Value = 1.539986
$ws.Range("I16").Value = 0.001874750541781658
$ws.Range("J16").Value = 0.001874750541781658
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 56.84506433333333
$ws.Range("N16").Value = 170.535193
$ws.Range("O16").Value = 0.2242663767030476
$ws.Range("P16").Value = 0.2242663767030477
$ws.Range("Q16").Value = 29.18020108081089
$ws.Range("R16").Value = 262.6218097272979
$ws.Range("S16").Value = 0.000420443511227448
$ws.Range("T16").Value = 0.000420443511227448

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.5133286666666667
$ws.Range("H17").Value = 1.539986
$ws.Range("I17").Value = 0.001874750541781658
$ws.Range("J17").Value = 0.001874750541781658
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 19.21324333333333
$ws.Range("N17").Value = 57.63973
$ws.Range("O17").Value = 0.07580050295684103
$ws.Range("P17").Value = 0.07580050295684104
$ws.Range("Q17").Value = 9.862708582642222
$ws.Range("R17").Value = 9.862708582642222
$ws.Range("S17").Value = 0.0001421070339856599
$ws.Range("T17").Value = 0.0001421070339856599

